$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSheet")

$ws.Range("A33").Value = "PDLServicePackage"
$ws.Range("B33").Value = "ADSL Test"

$ws.Range("J23").Select()
